$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column (D) keeps its text formatting so values
# like "1.00" / "577.70" are not reinterpreted as numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "64.151.69"
$ws.Range("E2").Value = "  +0.14%  "

$ws.Range("D3").Value = "2.761.80"
$ws.Range("E3").Value = "  +0.61%  "

$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.13%  "

$ws.Range("D5").Value = "577.70"
$ws.Range("E5").Value = "  -0.64%  "

$ws.Range("D6").Value = "159.36"
$ws.Range("E6").Value = "  +0.52%  "

$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  +0.01%  "

$ws.Range("D8").Value = "0.605"
$ws.Range("E8").Value = "  -2.90%  "

$ws.Range("D9").Value = "0.111"
$ws.Range("E9").Value = "  -1.56%  "

$ws.Range("D10").Value = "0.166"
$ws.Range("E10").Value = "  +4.14%  "

$ws.Range("B11").Value = "Toncoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D11").Value = "5.79"
$ws.Range("E11").Value = "  -14.88%  "

$ws.Range("B12").Value = "Cardano"
$ws.Range("C12").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D12").Value = "0.389"
$ws.Range("E12").Value = "  -1.07%  "

$ws.Range("D13").Value = "3.250.24"
$ws.Range("E13").Value = "  +0.52%  "

$ws.Range("D14").Value = "27.05"
$ws.Range("E14").Value = "  -1.99%  "

$ws.Range("D15").Value = "63.730.96"
$ws.Range("E15").Value = "  -0.38%  "

$ws.Range("D17").Value = "2.763.25"
$ws.Range("E17").Value = "  +0.20%  "

$ws.Range("D18").Value = "12.20"
$ws.Range("E18").Value = "  +0.18%  "

$ws.Range("D19").Value = "4.88"
$ws.Range("E19").Value = "  -1.49%  "

$ws.Range("D20").Value = "358.52"
$ws.Range("E20").Value = "  -1.55%  "

$ws.Range("D21").Value = "6.79"
$ws.Range("E21").Value = "  -3.07%  "

$ws.Range("D22").Value = "0.999"
$ws.Range("E22").Value = "  +0.42%  "

$ws.Range("D23").Value = "0.536"
$ws.Range("E23").Value = "  -0.97%  "

$ws.Range("D24").Value = "65.62"
$ws.Range("E24").Value = "  -1.80%  "

$ws.Range("E25").Value = "  -1.39%  "

$ws.Range("D26").Value = "8.62"
$ws.Range("E26").Value = "  -0.22%  "

$ws.Range("E27").Value = "  +0.00%  "

$ws.Range("D28").Value = "0.0₃0911"
$ws.Range("E28").Value = "  -0.79%  "

$ws.Range("D29").Value = "7.37"
$ws.Range("E29").Value = "  +1.35%  "

$ws.Range("E30").Value = "  -2.88%  "

$ws.Range("E31").Value = "  -1.64%  "

$ws.Range("D32").Value = "169.22"
$ws.Range("E32").Value = "  -2.53%  "

$ws.Range("B33").Value = "NEARProtocol"
$ws.Range("C33").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D33").Value = "4.98"
$ws.Range("E33").Value = "  +0.86%  "

$ws.Range("B34").Value = "EthereumClassic"
$ws.Range("C34").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D34").Value = "20.29"
$ws.Range("E34").Value = "  -1.75%  "

$ws.Range("E35").Value = "  +1.78%  "

$ws.Range("D36").Value = "0.998"
$ws.Range("E36").Value = "  -0.05%  "

$ws.Range("D37").Value = "1.82"
$ws.Range("E37").Value = "  -0.19%  "

$ws.Range("E38").Value = "  -0.57%  "

$ws.Range("B39").Value = "RenderToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D39").Value = "6.45"
$ws.Range("E39").Value = "  +4.69%  "

$ws.Range("B40").Value = "Bittensor"
$ws.Range("C40").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D40").Value = "348.35"
$ws.Range("E40").Value = "  +2.54%  "

$ws.Range("D41").Value = "4.21"
$ws.Range("E41").Value = "  -1.54%  "

$ws.Range("D42").Value = "39.17"
$ws.Range("E42").Value = "  -1.32%  "

$ws.Range("D43").Value = "21.55"
$ws.Range("E43").Value = "  -1.68%  "

$ws.Range("D44").Value = "21.91"
$ws.Range("E44").Value = "  -2.46%  "

$ws.Range("D45").Value = "0.0593"
$ws.Range("E45").Value = "  -1.93%  "

$ws.Range("B46").Value = "VeChain"
$ws.Range("C46").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D46").Value = "0.0256"
$ws.Range("E46").Value = "  -1.00%  "

$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D47").Value = "136.81"
$ws.Range("E47").Value = "  -0.77%  "

$ws.Range("B48").Value = "Mantle"
$ws.Range("C48").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D48").Value = "0.634"
$ws.Range("E48").Value = "  -1.77%  "

$ws.Range("E49").Value = "  -0.45%  "

$ws.Range("E50").Value = "  -0.13%  "

$ws.Range("E51").Value = "  -0.16%  "
